$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.260.08'
$ws.Range("E2").Value = '  +0.17%  '

$ws.Range("D3").Value = '1.857.10'
$ws.Range("E3").Value = '  -0.58%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.81'
$ws.Range("E5").Value = '  -2.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4560'
$ws.Range("E7").Value = '  -2.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3936'
$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.51'
$ws.Range("E9").Value = '  +0.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07808'
$ws.Range("E10").Value = '  -2.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9817'
$ws.Range("E11").Value = '  -2.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.28'
$ws.Range("E12").Value = '  -1.98%  '

$ws.Range("D13").Value = '1.851.83'
$ws.Range("E13").Value = '  -0.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.803'
$ws.Range("E14").Value = '  -3.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.958'
$ws.Range("E15").Value = '  -4.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.00'
$ws.Range("E17").Value = '  -3.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06528'
$ws.Range("E18").Value = '  -1.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001016'
$ws.Range("E19").Value = '  -2.38%  '

$ws.Range("E20").Value = '  -3.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").Value = '28.264.83'
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.284'
$ws.Range("E23").Value = '  -2.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.82'
$ws.Range("E24").Value = '  -2.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.252'
$ws.Range("E25").Value = '  -1.83%  '

$ws.Range("D26").Value = '2.072.86'
$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.91'
$ws.Range("E27").Value = '  -1.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.20'
$ws.Range("E28").Value = '  -4.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.055'
$ws.Range("E29").Value = '  -3.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.291'
$ws.Range("E30").Value = '  -3.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '116.39'
$ws.Range("E31").Value = '  -2.82%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9367'
$ws.Range("E32").Value = '  -3.83%  '

$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09267'
$ws.Range("E33").Value = '  -2.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.599'
$ws.Range("E34").Value = '  +0.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.377'
$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.195'
$ws.Range("E36").Value = '  -2.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06011'
$ws.Range("E37").Value = '  -1.47%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02210'
$ws.Range("E38").Value = '  -2.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.238'
$ws.Range("E39").Value = '  -2.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.158'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.0000'
$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5689'
$ws.Range("E42").Value = '  -4.69%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1792'
$ws.Range("E43").Value = '  -4.82%  '

$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.958'
$ws.Range("E44").Value = '  -3.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.258'
$ws.Range("E45").Value = '  -2.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.307'
$ws.Range("E46").Value = '  +17.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5392'
$ws.Range("E47").Value = '  -4.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.80'
$ws.Range("E48").Value = '  -2.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07156'
$ws.Range("E49").Value = '  +4.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.865'
$ws.Range("E50").Value = '  -5.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.77'
$ws.Range("E51").Value = '  -1.37%  '
